# This script applies a "sliding window" update to the sensor-reading
# dataset on Sheet1: four brand-new rows of accelerometer/gyroscope
# readings are inserted at the top of the data block (rows 2-5), the
# previously-existing readings in C2:H17 slide down to C6:H21, and the
# values that used to occupy the bottom of the window (old rows 18-21)
# fall out of the (fixed-size) window. Columns A (timestamp) and B
# (label) are left untouched, since the window itself doesn't move for
# those - only the feature columns C:H are resampled.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Capture the values currently sitting in C2:H17 - these are the
#    ones that need to slide down to C6:H21.
$shiftRange = $ws.Range("C2:H17")
$shiftedValues = $shiftRange.Value2

# 2. Write the brand new readings into the top of the window (rows 2-5).
$newTopValues = @(
    @(-0.7477993965148926, 0.8460922241210938, 0.0258597135543823, 0.01418807215633853, 0.06712245657330498, -0.07486735071454727),
    @(-0.8064756393432617, 0.8524413108825684, -0.0703473389148712, -0.02585268907603775, -0.07008951618557867, -0.06299911678901748),
    @(-0.7388706207275391, 0.7278079986572266, 0.0533058643341064, 0.046578474342823, -0.1310305893421173, -0.0218384321779012),
    @(-0.7223987579345703, 0.7231974601745605, 0.1768441945314407, 0.0740674127425465, -0.09423323614256736, -0.02838341776458984)
)

for ($i = 0; $i -lt 4; $i++) {
    $rowNum = 2 + $i
    for ($j = 0; $j -lt 6; $j++) {
        $colNum = 3 + $j
        $ws.Cells.Item($rowNum, $colNum).Value = $newTopValues[$i][$j]
    }
}

# 3. Write the previously-captured C2:H17 values down into C6:H21.
$ws.Range("C6:H21").Value2 = $shiftedValues
